$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 589
$ws.Range("I2").Value = 1.51025641025641
$ws.Range("J2").Value = 0.9914529914529915

$ws.Range("H3").Value = 1163
$ws.Range("I3").Value = 1.100283822138127
$ws.Range("J3").Value = 0.9395711500974658

$ws.Range("H4").Value = 138
$ws.Range("I4").Value = 1.112903225806452

$ws.Range("H5").Value = 455
$ws.Range("I5").Value = 1.350148367952522
$ws.Range("J5").Value = 0.9772502472799208

$ws.Range("H6").Value = 961
$ws.Range("I6").Value = 1.74410163339383
$ws.Range("J6").Value = 0.9685420447670902

$ws.Range("H7").Value = 287
$ws.Range("I7").Value = 2.87

$ws.Range("H8").Value = 104
$ws.Range("I8").Value = 2.536585365853659

$ws.Range("H9").Value = 263
$ws.Range("I9").Value = 1.91970802919708

$ws.Range("H10").Value = 1346
$ws.Range("I10").Value = 1.296724470134875
$ws.Range("J10").Value = 0.9473346178548491

$ws.Range("H12").Value = 206
$ws.Range("I12").Value = 2.102040816326531

$ws.Range("H13").Value = 1409
$ws.Range("I13").Value = 1.216753022452504
$ws.Range("J13").Value = 0.9242947610823259

$ws.Range("E14").Value = 1027
$ws.Range("F14").Value = 1.266337854500617
$ws.Range("G14").Value = 0.8947801068639539
$ws.Range("H14").Value = 1114
$ws.Range("I14").Value = 1.373612823674476
$ws.Range("J14").Value = 0.9905466502260584

$ws.Range("E15").Value = 151
$ws.Range("F15").Value = 2.126760563380282
$ws.Range("G15").Value = 0.7746478873239436
$ws.Range("H15").Value = 146
$ws.Range("I15").Value = 2.056338028169014

$ws.Range("H16").Value = 564
$ws.Range("I16").Value = 2.128301886792453
$ws.Range("J16").Value = 0.9974842767295597

$ws.Range("H17").Value = 3758
$ws.Range("I17").Value = 1.373538011695906
$ws.Range("J17").Value = 0.9522417153996101

$ws.Range("E18").Value = 4698
$ws.Range("F18").Value = 1.452690166975881
$ws.Range("G18").Value = 0.8593073593073594
$ws.Range("H18").Value = 4497
$ws.Range("I18").Value = 1.390538033395176
$ws.Range("J18").Value = 0.9706246134817563

$ws.Range("E19").Value = 507
$ws.Range("F19").Value = 1.198581560283688
$ws.Range("G19").Value = 0.8912529550827423
$ws.Range("H19").Value = 541
$ws.Range("I19").Value = 1.278959810874704
$ws.Range("J19").Value = 0.9976359338061466

$ws.Range("E20").Value = 527
$ws.Range("F20").Value = 1.93040293040293
$ws.Range("G20").Value = 0.9157509157509157
$ws.Range("H20").Value = 435
$ws.Range("I20").Value = 1.593406593406593

$ws.Range("H21").Value = 698
$ws.Range("I21").Value = 2.758893280632411
$ws.Range("J21").Value = 0.9960474308300395

$ws.Range("H23").Value = 126
$ws.Range("I23").Value = 3.230769230769231

$ws.Range("H24").Value = 134
$ws.Range("I24").Value = 5.36

$ws.Range("H25").Value = 1063
$ws.Range("I25").Value = 3.045845272206304
$ws.Range("J25").Value = 0.9923591212989494

$ws.Range("H26").Value = 298
$ws.Range("I26").Value = 4.082191780821918
$ws.Range("J26").Value = 0.9954337899543378

$ws.Range("E27").Value = 896
$ws.Range("F27").Value = 4.950276243093922
$ws.Range("G27").Value = 0.8066298342541437
$ws.Range("H27").Value = 831
$ws.Range("I27").Value = 4.591160220994476

$ws.Range("H28").Value = 317
$ws.Range("I28").Value = 1.42152466367713
$ws.Range("J28").Value = 0.9955156950672646

$ws.Range("H29").Value = 154
$ws.Range("I29").Value = 4.666666666666667

$ws.Range("H30").Value = 333
$ws.Range("I30").Value = 3.542553191489362
$ws.Range("J30").Value = 0.9964539007092199

$ws.Range("H31").Value = 44
$ws.Range("I31").Value = 1.466666666666667

$ws.Range("H32").Value = 1327
$ws.Range("I32").Value = 34.02564102564103

$ws.Range("H33").Value = 450
$ws.Range("I33").Value = 4.017857142857143

$ws.Range("E34").Value = 1559
$ws.Range("F34").Value = 16.58510638297872
$ws.Range("G34").Value = 0.8723404255319149
$ws.Range("H34").Value = 1538
$ws.Range("I34").Value = 16.36170212765957

$ws.Range("H35").Value = 954
$ws.Range("I35").Value = 7.511811023622047

$ws.Range("H36").Value = 893
$ws.Range("I36").Value = 30.79310344827586

$ws.Range("H37").Value = 819
$ws.Range("I37").Value = 15.45283018867925

$ws.Range("H38").Value = 1384
$ws.Range("I38").Value = 2.423817863397548
$ws.Range("J38").Value = 0.978984238178634

$ws.Range("H39").Value = 217
$ws.Range("I39").Value = 4.428571428571429

$ws.Range("H40").Value = 1330
$ws.Range("I40").Value = 13.43434343434343

$ws.Range("H41").Value = 307
$ws.Range("I41").Value = 12.79166666666667

$ws.Range("E42").Value = 1810
$ws.Range("F42").Value = 9.576719576719576
$ws.Range("G42").Value = 0.9435626102292769
$ws.Range("H42").Value = 1355
$ws.Range("I42").Value = 7.169312169312169

$ws.Range("H43").Value = 205
$ws.Range("I43").Value = 25.625

$ws.Range("H44").Value = 2456
$ws.Range("I44").Value = 79.2258064516129

$ws.Range("H45").Value = 1121
$ws.Range("I45").Value = 62.27777777777778
